$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.468.63"
$ws.Range("E2").Value = "  -1.94%  "
$ws.Range("D3").Value = "3.163.26"
$ws.Range("E3").Value = "  -3.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.28%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.164.43"
$ws.Range("E8").Value = "  -3.84%  "
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.96%  "
$ws.Range("E13").Value = "  -5.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.96%  "
$ws.Range("D15").Value = "3.666.25"
$ws.Range("E15").Value = "  -4.46%  "
$ws.Range("E16").Value = "  -2.09%  "
$ws.Range("D17").Value = "3.155.35"
$ws.Range("E17").Value = "  -4.14%  "
$ws.Range("D18").Value = "62.453.03"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("E19").Value = "  -4.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "452.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.43%  "
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.702"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.04%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.77%  "
$ws.Range("B28").Value = "FirstDigitalUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.104"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("E34").Value = "  -5.75%  "
$ws.Range("E35").Value = "  -6.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.83%  "
$ws.Range("D38").Value = "0.0₃0699"
$ws.Range("E38").Value = "  -5.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0384"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "403.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.62%  "
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("D44").Value = "2.794.83"
$ws.Range("E44").Value = "  -8.81%  "
$ws.Range("E45").Value = "  -5.33%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.111"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.29%  "
